# Correction génération données fichiers json lorsque une donnée de
# l'étudiant n'existe pas dans la base de données (ues, semestre ou
# années de formation exécutée par l'étudiant) -> remplacée par la
# valeur NULL.
#
# Concretely: column A (Numéro) values move from the 2015xxxx series
# to the 2017xxxx series, and several column E (Moyenne de l'étudiant)
# values are updated / regenerated for rows 3..63.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AI")

$changes = @(
    @{Row=3; A=20170926; E=8},
    @{Row=4; A=20170927; E=19},
    @{Row=5; A=20170928; E=$null},
    @{Row=6; A=20170929; E=20},
    @{Row=7; A=20170930; E=11},
    @{Row=8; A=20170931; E=$null},
    @{Row=9; A=20170932; E=7},
    @{Row=10; A=20170933; E=18},
    @{Row=11; A=20170934; E=11},
    @{Row=12; A=20170935; E=14},
    @{Row=13; A=20170936; E=6},
    @{Row=14; A=20170937; E=14},
    @{Row=15; A=20170938; E=11},
    @{Row=16; A=20170939; E=5},
    @{Row=17; A=20170940; E=9},
    @{Row=18; A=20170941; E=11},
    @{Row=19; A=20170942; E=20},
    @{Row=20; A=20170943; E=11},
    @{Row=21; A=20170944; E=17},
    @{Row=22; A=20170945; E=10},
    @{Row=23; A=20170946; E=$null},
    @{Row=24; A=20170947; E=10},
    @{Row=25; A=20170948; E=$null},
    @{Row=26; A=20170949; E=19},
    @{Row=27; A=20170950; E=17},
    @{Row=28; A=20170951; E=19},
    @{Row=29; A=20170952; E=$null},
    @{Row=30; A=20170953; E=17},
    @{Row=31; A=20170954; E=10},
    @{Row=32; A=20170955; E=9},
    @{Row=33; A=20170956; E=7},
    @{Row=34; A=20170957; E=9},
    @{Row=35; A=20170958; E=6},
    @{Row=36; A=20170959; E=16},
    @{Row=37; A=20170960; E=10},
    @{Row=38; A=20170961; E=5},
    @{Row=39; A=20170962; E=9},
    @{Row=40; A=20170963; E=5},
    @{Row=41; A=20170964; E=$null},
    @{Row=42; A=20170965; E=14},
    @{Row=43; A=20170966; E=17},
    @{Row=44; A=20170967; E=20},
    @{Row=45; A=20170968; E=10},
    @{Row=46; A=20170969; E=18},
    @{Row=47; A=20170970; E=$null},
    @{Row=48; A=20170971; E=11},
    @{Row=49; A=20170972; E=17},
    @{Row=50; A=20170973; E=8},
    @{Row=51; A=20170974; E=$null},
    @{Row=52; A=20170975; E=16},
    @{Row=53; A=20170976; E=8},
    @{Row=54; A=20170977; E=7},
    @{Row=55; A=20170978; E=15},
    @{Row=56; A=20170979; E=5},
    @{Row=57; A=20170980; E=15},
    @{Row=58; A=20170981; E=16},
    @{Row=59; A=20170982; E=17},
    @{Row=60; A=20170983; E=17},
    @{Row=61; A=20170984; E=6},
    @{Row=62; A=20170985; E=16},
    @{Row=63; A=20170986; E=10}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 1).Value = $change.A
    if ($null -ne $change.E) {
        $ws.Cells.Item($change.Row, 5).Value = $change.E
    }
}
